# Add the pull-request reference hyperlink for Annika Prasanna (row 18, column B)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B18 is currently empty (row for "Annika Prasanna" in column A).
# Insert a hyperlink whose display text is the PR URL itself - this is the
# same as Excel auto-linking a pasted/typed URL, which also creates the
# built-in "Hyperlink" cell style the first time it is used in the workbook.
$ws.Hyperlinks.Add($ws.Range("B18"), "https://github.com/dhavalkeerthi/MRIInterns2026A/pull/27") | Out-Null

# Leave the selection on the newly-linked cell, matching the edited workbook.
$ws.Range("B18").Select() | Out-Null
